$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 (secondary header row: "pompes)/Hiver/Eté/Année" style labels)
# is removed entirely; the data rows shift up by one.
$ws.Rows.Item(2).Delete() | Out-Null

# Rewrite the header row (row 1) with the new column headers.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# A1:E1 keep the plain/default cell format (no explicit style), even
# though some of them may have inherited a style from the deleted row.
$ws.Range("A1:E1").Style = "Normal"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 pick up a distinct style: same font as the rest of the sheet
# (Arial, 9pt) but with the default/general number format. Build this
# via a transient named style (so we reuse the existing Arial-9 font
# instead of minting a new one), then drop the named style again so we
# are left with a plain, unnamed cell format.
$tempStyle = $wb.Styles.Add("TempHeaderStyle")
$tempStyle.Font.Name = "Arial"
$tempStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TempHeaderStyle"
$tempStyle.Delete() | Out-Null

# Match the recorded selection state after the edit.
$ws.Range("A2:K2").Select() | Out-Null
